# Updates the "cryptos" price/volume list in place (Wed May 31 10:45:11 UTC 2023
# refresh). For each changed coin row we overwrite the Price (column D) and
# Volume(1h) (column E) text values; D46/D47 also got a coin swap
# (Decentraland <-> PaxDollar) in this run, so B/C/D/E are rewritten for
# those two rows.
#
# Price values are stored as plain text in the sheet (e.g. "27.113.28",
# "1.001", "0.4634") rather than numbers, so values that would otherwise be
# auto-parsed as a number by Excel (a single "." present) are entered with a
# leading apostrophe to force text entry and keep the exact original text
# (trailing zeros, etc.) intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.113.28"
$ws.Range("E2").Value = "  -2.99%  "
$ws.Range("D3").Value = "1.869.54"
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'307.55"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "'0.5057"
$ws.Range("E7").Value = "  +1.07%  "
$ws.Range("E8").Value = "  -1.89%  "
$ws.Range("D9").Value = "'0.07147"
$ws.Range("E9").Value = "  -2.37%  "
$ws.Range("D10").Value = "'0.8872"
$ws.Range("E10").Value = "  -2.87%  "
$ws.Range("D11").Value = "'20.64"
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("D12").Value = "'0.07554"
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").Value = "1.845.20"
$ws.Range("E13").Value = "  -3.29%  "
$ws.Range("D14").Value = "'5.318"
$ws.Range("E14").Value = "  -3.57%  "
$ws.Range("D15").Value = "'89.26"
$ws.Range("E15").Value = "  -3.68%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "'0.000008466"
$ws.Range("E17").Value = "  -3.25%  "
$ws.Range("D18").Value = "'14.10"
$ws.Range("E18").Value = "  -3.82%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").Value = "27.160.21"
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("D21").Value = "'5.072"
$ws.Range("E21").Value = "  -2.19%  "
$ws.Range("D22").Value = "2.095.41"
$ws.Range("E22").Value = "  -1.70%  "
$ws.Range("D23").Value = "'10.55"
$ws.Range("E23").Value = "  -2.82%  "
$ws.Range("D24").Value = "'6.474"
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").Value = "'150.92"
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("D26").Value = "'1.839"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").Value = "'18.00"
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("D28").Value = "'2.105"
$ws.Range("E28").Value = "  -5.42%  "
$ws.Range("D29").Value = "'112.71"
$ws.Range("E29").Value = "  -2.39%  "
$ws.Range("D30").Value = "'4.754"
$ws.Range("E30").Value = "  -3.33%  "
$ws.Range("D31").Value = "'4.689"
$ws.Range("E31").Value = "  -3.62%  "
$ws.Range("D32").Value = "'0.09038"
$ws.Range("E32").Value = "  +0.14%  "
$ws.Range("D33").Value = "'0.05125"
$ws.Range("E33").Value = "  -3.03%  "
$ws.Range("D34").Value = "'3.092"
$ws.Range("E34").Value = "  -3.61%  "
$ws.Range("D35").Value = "'0.7396"
$ws.Range("E35").Value = "  -4.73%  "
$ws.Range("D36").Value = "'1.158"
$ws.Range("E36").Value = "  -6.46%  "
$ws.Range("D37").Value = "'0.02033"
$ws.Range("E37").Value = "  -2.69%  "
$ws.Range("D38").Value = "'2.500"
$ws.Range("E38").Value = "  -3.59%  "
$ws.Range("D39").Value = "'3.043"
$ws.Range("E39").Value = "  -0.95%  "
$ws.Range("E40").Value = "  -1.41%  "
$ws.Range("D41").Value = "'0.5364"
$ws.Range("E41").Value = "  -3.65%  "
$ws.Range("D42").Value = "'6.591"
$ws.Range("E42").Value = "  -4.33%  "
$ws.Range("D43").Value = "'115.77"
$ws.Range("E43").Value = "  +2.62%  "
$ws.Range("D44").Value = "'8.416"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("D45").Value = "'0.1471"
$ws.Range("E45").Value = "  -3.41%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'1.001"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.4634"
$ws.Range("E47").Value = "  -4.27%  "
$ws.Range("D48").Value = "'10.01"
$ws.Range("E48").Value = "  -6.10%  "
$ws.Range("D49").Value = "'1.564"
$ws.Range("E49").Value = "  -4.71%  "
$ws.Range("D50").Value = "'64.49"
$ws.Range("E50").Value = "  -4.66%  "
$ws.Range("E51").Value = "  -1.68%  "
